# The "Add an image for a category (database change needed)." bullet
# item was highlighted in yellow as a to-do flag. The feature is now
# fully implemented, so remove the highlight from the whole paragraph
# (all of its runs, plus the paragraph mark itself).

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*category (database change needed)*") {
        $targetIndex = $i
    }
}

if ($targetIndex -gt 0) {
    $target = $d.Paragraphs.Item($targetIndex)
    $r = $target.Range
    $r.Font.HighlightColorIndex = 0   # wdNoHighlight
}
